# Zeiten.xlsx - "Zeiten angepasst und Codeformatierung"
# Add four new time-tracking rows (27-30) below the existing table and
# widen column E slightly to fit the new, longer "Tätigkeit" texts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append new rows -------------------------------------------------
# Copy the formatting (date/time number formats etc.) of the last
# existing row (26) down into the four new rows so the new cells reuse
# the same cell styles as the rest of the table.
[void]$ws.Range("A26:F26").Copy($ws.Range("A27:F30"))

# Row 27: 28.04.2015, 20:00 - 22:30
$ws.Range("A27").Value = 42122
$ws.Range("B27").Value = 0.83333333333333337
$ws.Range("C27").Value = 0.9375
$ws.Range("D27").Value = "Juliano,Tommy, Constanrin"
$ws.Range("E27").Value = "Einbau einer Fliesenklasse um das Programm Multiple-Device fähig zu machen"
$ws.Range("F27").ClearContents()

# Row 28: 30.04.2015, 23:00 - 00:10
$ws.Range("A28").Value = 42124
$ws.Range("B28").Value = 0.95833333333333337
$ws.Range("C28").Value = 0.006944444444444444
$ws.Range("D28").Value = "Juliano"
$ws.Range("E28").Value = "Schlange bewegt sich im Spielfeld"
$ws.Range("F28").ClearContents()

# Row 29: 02.05.2015, 20:00 - 03:30
$ws.Range("A29").Value = 42126
$ws.Range("B29").Value = 0.83333333333333337
$ws.Range("C29").Value = 0.14583333333333334
$ws.Range("D29").Value = "Juliano"
$ws.Range("E29").Value = "Schlange kann Beere fressen und wachsten, zudem ist das Spiel vorbei wenn Schlange mit sich selber kollidiert"
$ws.Range("F29").ClearContents()

# Row 30: 11.05.2015, 23:00 - 02:40
$ws.Range("A30").Value = 42135
$ws.Range("B30").Value = 0.95833333333333337
$ws.Range("C30").Value = 0.1111111111111111
$ws.Range("D30").Value = "Tommy"
$ws.Range("E30").Value = "Schwierigkeitsgrage mit Sound implementiert, BluetoothActivity erstellt"
$ws.Range("F30").ClearContents()

# --- Widen column E (Tätigkeit) so the longer texts fit --------------
$ws.Columns("E").ColumnWidth = 100.16666666666667

# --- Update selection / scroll position to match the edited view -----
[void]$excel.Goto($ws.Range("A10"), $true)
[void]$ws.Range("A31").Select()
